$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2) Model parameters")

# Fill in parameter descriptions (column B) that were previously blank,
# and correct the wording of an existing description (PropWetlands row).
# The order below matches the order the strings were (re)entered by the
# author, which determines the order they land in the shared string table.
$ws.Range("B19").Value = "Proportion of lake inflow as groundwater"
$ws.Range("B17").Value = "Loading rate of POC from wetlands"
$ws.Range("B16").Value = "Proportion of catchment that is wetlands"
$ws.Range("B8").Value  = "Lake water DOC concentration at start of model run"
$ws.Range("B9").Value  = "Lake water POC concentration at start of model run"
$ws.Range("B18").Value = "DOC concentration of groundwater"
$ws.Range("B20").Value = "DOC concentration of precipitation"
$ws.Range("B26").Value = "Influx of aerial POC (e.g., leaf litter)"

# Move the selection to reflect where the author ended up editing last (B26)
$ws.Range("B26").Select() | Out-Null
